# Weekly update: insert two new rows of data (week of 2023-12-07) right
# after the existing row 515, pushing all subsequent rows down by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 516/517 (everything below shifts down).
$ws.Rows.Item(516).Insert()
$ws.Rows.Item(517).Insert()

# New row 516: Betarraga, Primera
$ws.Cells.Item(516, 1).Value = 8
$ws.Cells.Item(516, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(516, 3).Value = "Coquimbo"
$ws.Cells.Item(516, 4).Value = 45267
$ws.Cells.Item(516, 5).Value = 4
$ws.Cells.Item(516, 6).Value = 100114014
$ws.Cells.Item(516, 7).Value = "Betarraga"
$ws.Cells.Item(516, 8).Value = "Sin especificar"
$ws.Cells.Item(516, 9).Value = "Primera"
$ws.Cells.Item(516, 10).Value = 1800
$ws.Cells.Item(516, 11).Value = 550
$ws.Cells.Item(516, 12).Value = 600
$ws.Cells.Item(516, 13).Value = 575
$ws.Cells.Item(516, 14).Value = "$/paquete 3 unidades"
$ws.Cells.Item(516, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(516, 16).Value = 192
$ws.Cells.Item(516, 17).Value = 3
$ws.Cells.Item(516, 18).Value = "Hortaliza"

# New row 517: Betarraga, Segunda
$ws.Cells.Item(517, 1).Value = 8
$ws.Cells.Item(517, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(517, 3).Value = "Coquimbo"
$ws.Cells.Item(517, 4).Value = 45267
$ws.Cells.Item(517, 5).Value = 4
$ws.Cells.Item(517, 6).Value = 100114014
$ws.Cells.Item(517, 7).Value = "Betarraga"
$ws.Cells.Item(517, 8).Value = "Sin especificar"
$ws.Cells.Item(517, 9).Value = "Segunda"
$ws.Cells.Item(517, 10).Value = 1000
$ws.Cells.Item(517, 11).Value = 450
$ws.Cells.Item(517, 12).Value = 500
$ws.Cells.Item(517, 13).Value = 475
$ws.Cells.Item(517, 14).Value = "$/paquete 3 unidades"
$ws.Cells.Item(517, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(517, 16).Value = 158
$ws.Cells.Item(517, 17).Value = 3
$ws.Cells.Item(517, 18).Value = "Hortaliza"
